$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column G width (matches width="17" in saved XML)
$ws.Columns.Item(7).ColumnWidth = 16.17

# Header cell G1 = "PRESUPUESTO", formatted like F1 (bold/centered header style)
$ws.Range("G1").Value = "PRESUPUESTO"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Data cells G2:G55 = 0, formatted like F2:F55 (currency style)
$ws.Range("G2:G55").Value = 0
$ws.Range("F2:F55").Copy()
$ws.Range("G2:G55").PasteSpecial(-4122)

# Totals row cell G56 = 0, formatted like F56 (totals currency style)
$ws.Range("G56").Value = 0
$ws.Range("F56").Copy()
$ws.Range("G56").PasteSpecial(-4122)

$excel.CutCopyMode = 0
